$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1).Range
$insertAt = $p1.End - 1

# Run 1: " ("
$r = $d.Range($insertAt, $insertAt)
$r.InsertAfter(" (")
$r.Font.Hidden = 1
$r.Font.Hidden = 0

# Run 2: "Changed main"
$insertAt = $r.End
$r = $d.Range($insertAt, $insertAt)
$r.InsertAfter("Changed main")
$r.Font.Hidden = 1
$r.Font.Hidden = 0

# Run 3: ")"
$insertAt = $r.End
$r = $d.Range($insertAt, $insertAt)
$r.InsertAfter(")")

Write-Output $d.Paragraphs(1).Range.Text
